$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 29; $r -ge 4; $r--) {
    $dst = $r + 2
    for ($col = 3; $col -le 20; $col++) {
        $ws.Cells.Item($dst, $col).Value2 = $ws.Cells.Item($r, $col).Value2
    }
}
Write-Host "done"
